# SwaadSutra_Consolidated_2026-01-21.xlsx update
# A new order (#28, Vipula Thakkar) was placed at 2026-01-21 10:01.
# It is inserted as the new first data row (row 2) of "All Orders",
# pushing every existing order down by one row, and the "Daily Summary"
# sheet's 2026-01-21 totals are updated accordingly.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "All Orders" ----
$ws = $wb.Worksheets.Item("All Orders")

# Make room for the new order at the top of the data (row 2).
$ws.Rows(2).Insert()

# New row 2 values
$ws.Range("A2").Value = 28
$ws.Range("B2").Value = "2026-01-21 10:01"
$ws.Range("C2").Value = "Vipula Thakkar"
$ws.Range("D2").Value = "B-903, Kakkad lavida"

# Phone and Collection Date look like numbers / dates but must stay text,
# matching the rest of the column.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "8109861246"

$ws.Range("F2").Value = "Appe Chutney x1"
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2026-01-21"

$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# ---- Sheet 2: "Daily Summary" ----
$ws2 = $wb.Worksheets.Item("Daily Summary")

# 2026-01-21 row: one more order (now 2), and the extra 60 revenue is
# still unpaid, so both Revenue and Pending increase by 60.
$ws2.Range("B2").Value = 2
$ws2.Range("E2").Value = 170
$ws2.Range("G2").Value = 170
